$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update country name positions (ranking shuffled due to updated case counts) ---
$ws.Range("A65").Value = "Costa Rica"
$ws.Range("A66").Value = "Moldavia"
$ws.Range("A78").Value = "Costa de Marfil"
$ws.Range("A79").Value = "Corea del Sur"
$ws.Range("A125").Value = "Ruanda"
$ws.Range("A126").Value = "Sri Lanka"
$ws.Range("A191").Value = "Brunei"
$ws.Range("A192").Value = "Seychelles"
$ws.Range("A193").Value = "Monaco"
$ws.Range("A202").Value = "Timor Oriental"
$ws.Range("A203").Value = "Santa Lucia"

# --- Update the "last updated" timestamp ---
$ws.Range("A1").Value = "Datos actualizados a 23 de Agosto de 2020 a las 23:21"

# --- Update numeric statistics per country row ---
$ws.Range("B4").Value = 5869507
$ws.Range("C4").Value = 28079
$ws.Range("D4").Value = 3160488
$ws.Range("E4").Value = 2528464
$ws.Range("G4").Value = 381
$ws.Range("H4").Value = 180555

$ws.Range("B5").Value = 3605783
$ws.Range("C5").Value = 23085
$ws.Range("E5").Value = 781401
$ws.Range("G5").Value = 467
$ws.Range("H5").Value = 114744

$ws.Range("B8").Value = 609773
$ws.Range("C8").Value = 2728
$ws.Range("D8").Value = 506470
$ws.Range("E8").Value = 90244
$ws.Range("G8").Value = 72
$ws.Range("H8").Value = 13059

$ws.Range("B27").Value = 124896
$ws.Range("C27").Value = 267
$ws.Range("D27").Value = 111112
$ws.Range("E27").Value = 4711

$ws.Range("B33").Value = 102663
$ws.Range("C33").Value = 730
$ws.Range("D33").Value = 79514
$ws.Range("E33").Value = 22315

$ws.Range("B35").Value = 91161
$ws.Range("C35").Value = 600
$ws.Range("D35").Value = 60670
$ws.Range("E35").Value = 28924
$ws.Range("G35").Value = 13
$ws.Range("H35").Value = 1567

$ws.Range("B65").Value = 33820
$ws.Range("C65").Value = 736
$ws.Range("D65").Value = 10518
$ws.Range("E65").Value = 22947
$ws.Range("G65").Value = 7
$ws.Range("H65").Value = 355

$ws.Range("B66").Value = 33478
$ws.Range("C66").Value = 406
$ws.Range("D66").Value = 22683
$ws.Range("E66").Value = 9855
$ws.Range("G66").Value = 5
$ws.Range("H66").Value = 940

$ws.Range("B78").Value = 17471
$ws.Range("C78").Value = 97
$ws.Range("D78").Value = 15301
$ws.Range("E78").Value = 2057
$ws.Range("H78").Value = 113

$ws.Range("B79").Value = 17399
$ws.Range("C79").Value = 397
$ws.Range("D79").Value = 14200
$ws.Range("E79").Value = 2890
$ws.Range("H79").Value = 309

$ws.Range("B90").Value = 10318
$ws.Range("C90").Value = 19
$ws.Range("E90").Value = 904

$ws.Range("B103").Value = 6905
$ws.Range("C103").Value = 11
$ws.Range("D103").Value = 6232
$ws.Range("E103").Value = 515

$ws.Range("B106").Value = 5854
$ws.Range("C106").Value = 316
$ws.Range("D106").Value = 2509
$ws.Range("E106").Value = 3293
$ws.Range("G106").Value = 6
$ws.Range("H106").Value = 52

$ws.Range("B125").Value = 3089
$ws.Range("C125").Value = 200
$ws.Range("D125").Value = 1755
$ws.Range("E125").Value = 1322

$ws.Range("B126").Value = 2953
$ws.Range("C126").Value = 6
$ws.Range("D126").Value = 2805
$ws.Range("E126").Value = 136
$ws.Range("G126").Value = 1
$ws.Range("H126").Value = 12

$ws.Range("B127").Value = 2818
$ws.Range("C127").Value = 80
$ws.Range("D127").Value = 1443
$ws.Range("E127").Value = 1304
$ws.Range("G127").Value = 3
$ws.Range("H127").Value = 71

$ws.Range("B191").Value = 143
$ws.Range("D191").Value = 139
$ws.Range("E191").Value = 1
$ws.Range("H191").Value = 3

$ws.Range("B192").Value = 132
$ws.Range("D192").Value = 126
$ws.Range("E192").Value = 6
$ws.Range("H192").Value = 0

$ws.Range("B193").Value = 112
$ws.Range("C193").Value = 4
$ws.Range("D193").Value = 83
$ws.Range("E193").Value = 28
$ws.Range("H193").Value = 1

